$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.068.81"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.16"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.23"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07488"
$ws.Range("E8").Value = "  -1.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2920"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.24"
$ws.Range("E10").Value = "  +2.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07675"
$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.837.39"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.003"
$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6671"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.67"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009333"
$ws.Range("E16").Value = "  -8.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.978"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.084.39"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.080.27"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "222.97"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.094"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.90"
$ws.Range("E25").Value = "  +0.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1392"
$ws.Range("E26").Value = "  +1.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.481"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.86"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05672"
$ws.Range("E30").Value = "  +8.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.146"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.072"
$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.207"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7411"
$ws.Range("E34").Value = "  +0.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.831"
$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.668"
$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.758"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.214.50"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01775"
$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.503"
$ws.Range("E41").Value = "  +2.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8888"
$ws.Range("E42").Value = "  -1.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.01"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.980.51"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("E46").Value = "  +0.74%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.59"
$ws.Range("E47").Value = "  +2.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5089"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4061"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07445"
$ws.Range("E50").Value = "  +6.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.981"
$ws.Range("E51").Value = "  +0.89%  "
